$d = $word.ActiveDocument

# 1. Remove the "Make sure to select your county for full results." paragraph
#    (including its own paragraph mark) to clean up the page.
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Make sure to select your county for full results.*") {
        $p.Range.Delete()
        $found = $true
        break
    }
}

if (-not $found) {
    # Fallback: locate the sentence with Find and expand to its paragraph.
    $rng = $d.Content
    $rng.Find.Execute("Make sure to select your county for full results.", $true,
                       $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($rng.Find.Found) {
        $rng.Expand(4) | Out-Null  # wdParagraph
        $rng.Delete()
    }
}

# 2. Reposition the "What's Next Text" graphic (anchored picture) that now
#    moves up/right to take the freed-up space.
$target = $null
for ($i = 1; $i -le $d.Shapes.Count; $i++) {
    $s = $d.Shapes.Item($i)
    if ($s.Name -like "Graphic*") {
        $target = $s
        break
    }
}
if (-not $target -and $d.Shapes.Count -ge 1) {
    $target = $d.Shapes.Item(1)
}
if ($target) {
    $target.Left = 391.9
    $target.Top = 260.4
}

# 3. The Heading 1 style (used by the "Your forms are ready. Well done!"
#    title) loses its explicit bold formatting.
$h1 = $d.Styles.Item("Heading 1")
$h1.Font.Bold = $false
$h1.Font.BoldBi = $false

$h1c = $d.Styles.Item("Heading 1 Char")
$h1c.Font.Bold = $false
$h1c.Font.BoldBi = $false
